# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1) Metadata sheet: bump the generation "Date" value.
# 2) Elements sheet: the two custom "Mapping:" columns (RIM Mapping / the
#    new French business-spec mapping) swap places - the new mapping
#    column (previously AL) moves to AK, and RIM Mapping (previously AK)
#    moves to AL. Header text, all data rows and the bestFit column
#    widths move together.

$wb = $excel.ActiveWorkbook

# --- 1) Metadata!B8 : Date -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2) Elements sheet: swap columns AK (37) and AL (38) ------------------
$ws = $wb.Worksheets.Item("Elements")

$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 1) { $lastRow = 6 }

for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $ws.Cells.Item($r, 37)
    $alCell = $ws.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    # Only touch cells whose content actually changes - an untouched round
    # trip keeps e.g. "typed empty string" cells exactly as they were,
    # whereas writing "" back through .Value collapses them to a plain
    # blank cell.
    if ($akVal -ne $alVal) {
        $akCell.Value = $alVal
        $alCell.Value = $akVal
    }
}

# Swap the (bestFit) column widths along with the content.
$ws.Columns.Item(37).ColumnWidth = 87
$ws.Columns.Item(38).ColumnWidth = 24.166666666666668
